# chartink_screener.xlsx - "break out stock.yaml completed"
#
# 1. On sheet "10per change": convert the bsecode column (E) for the
#    existing rows 4-6 from text-stored numbers to real numeric cells,
#    and append two new scraped rows (7 and 8) for the 06:45:35 run.
# 2. On sheet "DND 3 V 0.3": convert E2 (bsecode) from text to numeric.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 10per change
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("10per change")

# Existing rows: bsecode was stored as text ("541729" etc) - fix to numeric
$ws.Range("E4").Value = 541729
$ws.Range("E5").Value = 543237
$ws.Range("E6").Value = 526371

# New row 7 - MAZDOCK pulled at 06:45:35 (bsecode kept as text, like it was
# originally scraped)
$ws.Range("A7").Value = "10/06/2024 06:45:35"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "MAZDOCK"
$ws.Range("D7").Value = "Mazagon Dock Shipbuilders Ltd"
$ws.Range("E7").Value = "'543237"
$ws.Range("F7").Value = -1.09
$ws.Range("G7").Value = 3118.85
$ws.Range("H7").Value = 840927

# New row 8 - NMDC pulled at 06:45:35 (bsecode kept as text)
$ws.Range("A8").Value = "10/06/2024 06:45:35"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "NMDC"
$ws.Range("D8").Value = "Nmdc Limited"
$ws.Range("E8").Value = "'526371"
$ws.Range("F8").Value = -0.77
$ws.Range("G8").Value = 256.5
$ws.Range("H8").Value = 3161735

# ---------------------------------------------------------------------
# Sheet: DND 3 V 0.3
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DND 3 V 0.3")
$ws2.Range("E2").Value = 540755
